$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-87 down to 21-88.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record.
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 44910
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100112012
$ws.Cells.Item(20, 7).Value = "Espinaca"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 2700
$ws.Cells.Item(20, 12).Value = 3000
$ws.Cells.Item(20, 13).Value = 2850
$ws.Cells.Item(20, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 950
$ws.Cells.Item(20, 17).Value = 3
$ws.Cells.Item(20, 18).Value = "Hortaliza"
